# Add a new login-data row ("iameng" / "password") to the LoginTestData
# sheet, right after the existing iamfd/iamsup rows, then leave the
# selection parked on B7 (matching the post-edit state captured in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "iameng"
$ws.Range("B4").Value = "password"

$ws.Range("B7").Select()
